$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.503.55'
$ws.Range('E2').Value = '  -4.34%  '
$ws.Range('D3').Value = '3.356.06'
$ws.Range('E3').Value = '  -5.19%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.24'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.52%  '
$ws.Range('E7').Value = '  -2.18%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '3.347.57'
$ws.Range('E9').Value = '  -5.08%  '
$ws.Range('E10').Value = '  -8.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.592'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.67'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -8.14%  '
$ws.Range('E13').Value = '  -6.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.33%  '
$ws.Range('D15').Value = '3.889.12'
$ws.Range('E15').Value = '  -5.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '605.09'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -9.48%  '
$ws.Range('D17').Value = '66.556.66'
$ws.Range('E17').Value = '  -4.24%  '
$ws.Range('D18').Value = '3.360.22'
$ws.Range('E18').Value = '  -4.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.01'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.88%  '
$ws.Range('E20').Value = '  -3.25%  '
$ws.Range('E21').Value = '  -7.81%  '
$ws.Range('E22').Value = '  -5.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.06'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('E28').Value = '  -7.90%  '
$ws.Range('E29').Value = '  -8.20%  '
$ws.Range('E30').Value = '  -9.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.46'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.27'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.79'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -15.04%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.08'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.88%  '
$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '558.29'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +10.37%  '
$ws.Range('D36').Value = '3.829.36'
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('E37').Value = '  -5.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '58.14'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.15%  '
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('E40').Value = '  -7.67%  '
$ws.Range('D41').Value = '0.0₃0715'
$ws.Range('E41').Value = '  -12.20%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.127'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.58%  '
$ws.Range('B43').Value = 'CoreDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +22.60%  '
$ws.Range('E44').Value = '  -9.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.345'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '32.09'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0414'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.16'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.64'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.60%  '
$ws.Range('E50').Value = '  -4.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.05%  '

Write-Output "Applied changes"
